$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "C:\\R_S\\Instr\\user\\NR5G\\AllocationFiles\\DL\\64QAM_cellId1_papr11_74.allocation"
